# Updates the cryptos price/volume table with refreshed values.
# Numeric-looking "Price" strings are prefixed with a leading apostrophe so
# Excel stores them as text (preserving trailing zeros / dotted-thousands
# formatting) instead of silently coercing them into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.882.97"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "3.416.59"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'410.17"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "'128.95"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D7").Value = "'0.630"
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "'0.730"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").Value = "'0.138"
$ws.Range("E10").Value = "  -2.95%  "
$ws.Range("D11").Value = "'43.25"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "'9.16"
$ws.Range("E12").Value = "  +2.87%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.959.95"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "'0.141"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000210"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "'20.96"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").Value = "3.416.73"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "'1.09"
$ws.Range("E18").Value = "  +1.89%  "
$ws.Range("D19").Value = "'12.38"
$ws.Range("E19").Value = "  -1.75%  "
$ws.Range("D20").Value = "61.804.26"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").Value = "'484.69"
$ws.Range("E21").Value = "  +19.63%  "
$ws.Range("D22").Value = "'92.03"
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").Value = "'3.29"
$ws.Range("E23").Value = "  +3.22%  "
$ws.Range("D24").Value = "'13.48"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").Value = "'3.34"
$ws.Range("E25").Value = "  +3.60%  "
$ws.Range("D26").Value = "'34.04"
$ws.Range("E26").Value = "  +3.16%  "
$ws.Range("D27").Value = "'9.10"
$ws.Range("E27").Value = "  +5.91%  "
$ws.Range("D28").Value = "'4.80"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "'7.69"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").Value = "'2.77"
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("D31").Value = "'11.98"
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("D32").Value = "'0.114"
$ws.Range("E32").Value = "  -4.02%  "
$ws.Range("D33").Value = "'0.167"
$ws.Range("E33").Value = "  -2.58%  "
$ws.Range("D34").Value = "'41.77"
$ws.Range("E34").Value = "  -4.79%  "
$ws.Range("D36").Value = "'57.82"
$ws.Range("E36").Value = "  +9.14%  "
$ws.Range("D37").Value = "'0.0493"
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'150.51"
$ws.Range("E39").Value = "  +6.65%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'3.41"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "'0.136"
$ws.Range("E41").Value = "  +3.83%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "'0.324"
$ws.Range("E42").Value = "  +3.09%  "
$ws.Range("D43").Value = "'2.95"
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("D44").Value = "'2.12"
$ws.Range("E44").Value = "  +7.13%  "
$ws.Range("D45").Value = "'2.63"
$ws.Range("E45").Value = "  +10.66%  "
$ws.Range("D46").Value = "'4.25"
$ws.Range("E46").Value = "  +5.73%  "
$ws.Range("D47").Value = "'16.54"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("D48").Value = "'2.32"
$ws.Range("E48").Value = "  +19.38%  "
$ws.Range("D49").Value = "'22.56"
$ws.Range("E49").Value = "  +2.71%  "
$ws.Range("D50").Value = "'116.33"
$ws.Range("E50").Value = "  +21.49%  "
$ws.Range("E51").Value = "  +14.74%  "
